# Apply the "Rename keys, add tests, edit changelog, update spreadsheets" edit
# to the INTENT sheet of the Voxa Cli Intents/Utterances workbook:
#   - add two new columns "parameterName" / "parameterValue" (S1/T1)
#   - add a new intent row (row 9): HUMANINTENT / skill / Fallback skill
#   - tweak the "BOOLEAN" custom number format used by the slot/confirmation
#     boolean columns so the embedded literal "E" is lower-cased to "e"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INTENT")

# --- new header cells, matching the style of the existing header row ---
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null
$ws.Range("S1").Value = "parameterName"

$ws.Range("R1").Copy() | Out-Null
$ws.Range("T1").PasteSpecial(-4122) | Out-Null
$ws.Range("T1").Value = "parameterValue"

# --- new data row for the HUMANINTENT / Fallback skill handler ---
$ws.Range("A9").Value = "HUMANINTENT"
$ws.Range("S9").Value = "skill"
$ws.Range("T9").Value = "Fallback skill"

# --- re-case the custom "BOOLEAN" number format used by the boolean columns ---
foreach ($addr in @("H5", "P5", "Q5", "R5", "H6", "Q6", "R6", "H7", "Q7", "R7")) {
    $ws.Range($addr).NumberFormat = '"BOOL"e"AN"'
}
